$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 602
$ws.Range("B3").Value = "edrichhans"
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = "Adrian Sing"
$ws.Range("E3").Value = 3000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 60
$ws.Range("Q3").Value = 0
$ws.Range("S3").Value = -1750
